$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.726.13"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.23%  '

$ws.Range("D3").Value = "'2.976.18"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.72%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").Value = "'560.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.78%  '

$ws.Range("D6").Value = "'136.61"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +11.30%  '

$ws.Range("E7").Value = '  -0.09%  '

$ws.Range("D8").Value = "'0.517"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +4.44%  '

$ws.Range("D9").Value = "'2.970.39"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.60%  '

$ws.Range("D10").Value = "'0.131"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +6.27%  '

$ws.Range("D11").Value = "'4.85"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.10%  '

$ws.Range("D12").Value = "'0.455"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.04%  '

$ws.Range("D13").Value = "'0.0000226"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +7.33%  '

$ws.Range("D14").Value = "'33.44"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.38%  '

$ws.Range("E15").Value = '  +2.87%  '

$ws.Range("D16").Value = "'3.461.21"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.59%  '

$ws.Range("D17").Value = "'6.94"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +5.80%  '

$ws.Range("D18").Value = "'2.969.23"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.63%  '

$ws.Range("D19").Value = "'58.719.54"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.25%  '

$ws.Range("D20").Value = "'421.19"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.20%  '

$ws.Range("D21").Value = "'13.48"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.67%  '

$ws.Range("D22").Value = "'0.710"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +6.22%  '

$ws.Range("D23").Value = "'7.09"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.13%  '

$ws.Range("D24").Value = "'13.32"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.25%  '

$ws.Range("D25").Value = "'80.07"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.18%  '

$ws.Range("E26").Value = '  -0.06%  '

$ws.Range("E27").Value = '  +0.00%  '

$ws.Range("D28").Value = "'2.11"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +9.08%  '

$ws.Range("D29").Value = "'2.51"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.74%  '

$ws.Range("D30").Value = "'7.70"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +7.57%  '

$ws.Range("D31").Value = "'25.54"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.61%  '

$ws.Range("D32").Value = "'6.06"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.61%  '

$ws.Range("D33").Value = "'0.0989"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.64%  '

$ws.Range("D34").Value = "'1.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +10.75%  '

$ws.Range("D35").Value = "'0.0₃0762"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +23.44%  '

$ws.Range("D36").Value = "'5.71"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.59%  '

$ws.Range("D37").Value = "'2.06"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.30%  '

$ws.Range("D38").Value = "'48.79"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.83%  '

$ws.Range("D39").Value = "'8.64"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.50%  '

$ws.Range("D40").Value = "'2.74"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +14.33%  '

$ws.Range("D41").Value = "'397.27"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +10.34%  '

$ws.Range("D42").Value = "'2.740.90"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.72%  '

$ws.Range("D43").Value = "'0.0347"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.89%  '

$ws.Range("E44").Value = '  +0.74%  '

$ws.Range("D46").Value = "'124.87"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.93%  '

$ws.Range("D47").Value = "'0.243"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +6.20%  '

$ws.Range("D48").Value = "'2.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.30%  '

$ws.Range("E49").Value = '  +2.20%  '

$ws.Range("D50").Value = "'32.13"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +20.61%  '

$ws.Range("D51").Value = "'23.17"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.39%  '
